# Auto-generated script applying numeric corrections to Phoenix_Profits-style
# leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Re-derives currentAveragePrice / NQ / HQ price & profit columns (H:N) for a
# set of rows per sheet, matching the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 335258.34
$ws.Range("I43").Value = 401710
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 401710
$ws.Range("L43").Value = 3000
$ws.Range("M43").Value = -401641
$ws.Range("N43").Value = -3138
$ws.Range("H118").Value = 869.625
$ws.Range("I118").Value = 474.66666
$ws.Range("J118").Value = 2054.5
$ws.Range("K118").Value = 1423.99998
$ws.Range("L118").Value = 6163.5
$ws.Range("M118").Value = 233.0000199999999
$ws.Range("N118").Value = -9477.5
$ws.Range("H132").Value = 2766.4119
$ws.Range("I132").Value = 2746.1904
$ws.Range("K132").Value = 8238.5712
$ws.Range("M132").Value = -5708.5712
$ws.Range("H138").Value = 4999
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4999
$ws.Range("K138").Value = 0
$ws.Range("M138").Value = 14997
$ws.Range("N138").Value = -25277
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 39032.5
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 39032.5
$ws.Range("K37").Value = 0
$ws.Range("M37").Value = 39032.5
$ws.Range("N37").Value = -39578.5
$ws.Range("H45").Value = 2374.647
$ws.Range("I45").Value = 1170.1666
$ws.Range("J45").Value = 5265.4
$ws.Range("K45").Value = 1170.1666
$ws.Range("L45").Value = 5265.4
$ws.Range("M45").Value = -793.1666
$ws.Range("N45").Value = -6019.4
$ws.Range("H74").Value = 1717.9615
$ws.Range("I74").Value = 1736
$ws.Range("J74").Value = 1693.3636
$ws.Range("K74").Value = 1736
$ws.Range("L74").Value = 1693.3636
$ws.Range("M74").Value = -862
$ws.Range("N74").Value = -3441.3636
$ws.Range("H77").Value = 1717.9615
$ws.Range("I77").Value = 1736
$ws.Range("J77").Value = 1693.3636
$ws.Range("K77").Value = 8680
$ws.Range("L77").Value = 8466.817999999999
$ws.Range("M77").Value = -4312
$ws.Range("N77").Value = -17202.818
$ws.Range("H122").Value = 2254.8
$ws.Range("I122").Value = 2254.8
$ws.Range("K122").Value = 6764.400000000001
$ws.Range("M122").Value = -4314.400000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2436.4
$ws.Range("I86").Value = 2336.5715
$ws.Range("K86").Value = 2336.5715
$ws.Range("M86").Value = -1213.5715
$ws.Range("H89").Value = 2436.4
$ws.Range("I89").Value = 2336.5715
$ws.Range("K89").Value = 11682.8575
$ws.Range("M89").Value = -6066.8575
$ws.Range("H134").Value = 4274.6934
$ws.Range("I134").Value = 3633.8572
$ws.Range("J134").Value = 10255.833
$ws.Range("K134").Value = 10901.5716
$ws.Range("L134").Value = 30767.499
$ws.Range("M134").Value = -8366.571599999999
$ws.Range("N134").Value = -35837.499
$ws.Range("H138").Value = 85741
$ws.Range("J138").Value = 85741
$ws.Range("L138").Value = 85741
$ws.Range("N138").Value = -96021
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10368.685
$ws.Range("I31").Value = 3772.182
$ws.Range("J31").Value = 13056.148
$ws.Range("K31").Value = 3772.182
$ws.Range("L31").Value = 13056.148
$ws.Range("M31").Value = -3477.182
$ws.Range("N31").Value = -13646.148
$ws.Range("H34").Value = 10368.685
$ws.Range("I34").Value = 3772.182
$ws.Range("J34").Value = 13056.148
$ws.Range("K34").Value = 3772.182
$ws.Range("L34").Value = 13056.148
$ws.Range("M34").Value = -3570.182
$ws.Range("N34").Value = -13460.148
$ws.Range("H55").Value = 15000
$ws.Range("I55").Value = 5000
$ws.Range("J55").Value = 25000
$ws.Range("K55").Value = 5000
$ws.Range("L55").Value = 25000
$ws.Range("M55").Value = -4685
$ws.Range("N55").Value = -25630
$ws.Range("H58").Value = 3735.44
$ws.Range("I58").Value = 3704.1904
$ws.Range("J58").Value = 3899.5
$ws.Range("K58").Value = 3704.1904
$ws.Range("L58").Value = 3899.5
$ws.Range("M58").Value = -3501.1904
$ws.Range("N58").Value = -4305.5
$ws.Range("H59").Value = 58715.11
$ws.Range("J59").Value = 64179.5
$ws.Range("L59").Value = 64179.5
$ws.Range("N59").Value = -66469.5
$ws.Range("H60").Value = 35092.9
$ws.Range("I60").Value = 93
$ws.Range("J60").Value = 38981.777
$ws.Range("K60").Value = 93
$ws.Range("L60").Value = 38981.777
$ws.Range("M60").Value = 418
$ws.Range("N60").Value = -40003.777
$ws.Range("H132").Value = 8072.1943
$ws.Range("I132").Value = 6276.7334
$ws.Range("J132").Value = 17049.5
$ws.Range("K132").Value = 18830.2002
$ws.Range("L132").Value = 51148.5
$ws.Range("M132").Value = -16300.2002
$ws.Range("N132").Value = -56208.5
$ws.Range("H135").Value = 68866.664
$ws.Range("J135").Value = 68866.664
$ws.Range("L135").Value = 68866.664
$ws.Range("N135").Value = -79006.664
$ws.Range("H136").Value = 3735.44
$ws.Range("I136").Value = 3704.1904
$ws.Range("J136").Value = 3899.5
$ws.Range("K136").Value = 11112.5712
$ws.Range("L136").Value = 11698.5
$ws.Range("M136").Value = -8562.5712
$ws.Range("N136").Value = -16798.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1374.1333
$ws.Range("I33").Value = 90.833336
$ws.Range("J33").Value = 2229.6667
$ws.Range("K33").Value = 545.000016
$ws.Range("L33").Value = 13378.0002
$ws.Range("M33").Value = -262.000016
$ws.Range("N33").Value = -13944.0002
$ws.Range("H101").Value = 29850
$ws.Range("J101").Value = 29850
$ws.Range("L101").Value = 89550
$ws.Range("N101").Value = -94418
$ws.Range("H113").Value = 27779670
$ws.Range("I113").Value = 130.71428
$ws.Range("J113").Value = 34485076
$ws.Range("K113").Value = 392.14284
$ws.Range("L113").Value = 103455228
$ws.Range("M113").Value = 1777.85716
$ws.Range("N113").Value = -103459568
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 41021
$ws.Range("I26").Value = 37000
$ws.Range("K26").Value = 37000
$ws.Range("M26").Value = -36720
$ws.Range("H50").Value = 41021
$ws.Range("I50").Value = 37000
$ws.Range("K50").Value = 37000
$ws.Range("M50").Value = -36502
$ws.Range("H135").Value = 93268.17999999999
$ws.Range("J135").Value = 93268.17999999999
$ws.Range("L135").Value = 93268.17999999999
$ws.Range("N135").Value = -103408.18
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2949.1738
$ws.Range("J100").Value = 2822
$ws.Range("L100").Value = 2822
$ws.Range("N100").Value = -3904
$ws.Range("H122").Value = 5707.952
$ws.Range("I122").Value = 3409.3333
$ws.Range("K122").Value = 10227.9999
$ws.Range("M122").Value = -7777.999899999999
$ws.Range("H132").Value = 13816.158
$ws.Range("I132").Value = 9834.467000000001
$ws.Range("K132").Value = 29503.401
$ws.Range("M132").Value = -26973.401
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2086.9583
$ws.Range("I136").Value = 1908.9524
$ws.Range("J136").Value = 3333
$ws.Range("K136").Value = 5726.857199999999
$ws.Range("L136").Value = 9999
$ws.Range("M136").Value = -3176.857199999999
$ws.Range("N136").Value = -15099
$ws.Range("H139").Value = 74484.91
$ws.Range("J139").Value = 74484.91
$ws.Range("L139").Value = 74484.91
$ws.Range("N139").Value = -84764.91

# Columns that moved from NQ-profit (M) to HQ-profit (N) on these rows
# (or vice versa) -- clear the now-unused cell so the row shape matches.
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M138").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M37").ClearContents()

Write-Host "Applied Phoenix_Profits price/profit updates"
